{"js": "// Replace each three-digit x one-digit multiplication fact with its updated value.\n// Each 'find' string is unique within the document, so a direct search+replace\n// on the whole body is safe and order-independent.\nconst replacements = [\n  { find: \"380\u00d72=760\", replace: \"314\u00d79=2826\" },\n  { find: \"779\u00d74=3116\", replace: \"324\u00d79=2916\" },\n  { find: \"477\u00d73=1431\", replace: \"817\u00d75=4085\" },\n  { find: \"224\u00d72=448\", replace: \"151\u00d72=302\" },\n  { find: \"922\u00d72=1844\", replace: \"838\u00d77=5866\" },\n  { find: \"584\u00d76=3504\", replace: \"626\u00d76=3756\" },\n  { find: \"136\u00d77=952\", replace: \"936\u00d77=6552\" },\n  { find: \"490\u00d72=980\", replace: \"522\u00d73=1566\" },\n  { find: \"295\u00d72=590\", replace: \"740\u00d77=5180\" },\n  { find: \"962\u00d73=2886\", replace: \"379\u00d79=3411\" },\n  { find: \"623\u00d78=4984\", replace: \"274\u00d72=548\" },\n  { find: \"895\u00d73=2685\", replace: \"451\u00d78=3608\" },\n  { find: \"356\u00d75=1780\", replace: \"290\u00d77=2030\" },\n  { find: \"123\u00d76=738\", replace: \"983\u00d75=4915\" },\n  { find: \"519\u00d73=1557\", replace: \"530\u00d72=1060\" },\n  { find: \"768\u00d75=3840\", replace: \"181\u00d79=1629\" },\n  { find: \"989\u00d74=3956\", replace: \"280\u00d73=840\" },\n  { find: \"843\u00d77=5901\", replace: \"219\u00d74=876\" },\n  { find: \"512\u00d77=3584\", replace: \"415\u00d75=2075\" },\n  { find: \"595\u00d73=1785\", replace: \"553\u00d74=2212\" },\n  { find: \"571\u00d78=4568\", replace: \"307\u00d75=1535\" },\n  { find: \"797\u00d77=5579\", replace: \"471\u00d72=942\" },\n  { find: \"321\u00d79=2889\", replace: \"779\u00d78=6232\" },\n  { find: \"402\u00d76=2412\", replace: \"271\u00d75=1355\" },\n  { find: \"222\u00d76=1332\", replace: \"578\u00d75=2890\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication fact with its updated value.\n# Each 'find' string is unique within the document, so Find/Replace against the\n# whole document Content range is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = '380\u00d72=760'; Replace = '314\u00d79=2826' }\n    @{ Find = '779\u00d74=3116'; Replace = '324\u00d79=2916' }\n    @{ Find = '477\u00d73=1431'; Replace = '817\u00d75=4085' }\n    @{ Find = '224\u00d72=448'; Replace = '151\u00d72=302' }\n    @{ Find = '922\u00d72=1844'; Replace = '838\u00d77=5866' }\n    @{ Find = '584\u00d76=3504'; Replace = '626\u00d76=3756' }\n    @{ Find = '136\u00d77=952'; Replace = '936\u00d77=6552' }\n    @{ Find = '490\u00d72=980'; Replace = '522\u00d73=1566' }\n    @{ Find = '295\u00d72=590'; Replace = '740\u00d77=5180' }\n    @{ Find = '962\u00d73=2886'; Replace = '379\u00d79=3411' }\n    @{ Find = '623\u00d78=4984'; Replace = '274\u00d72=548' }\n    @{ Find = '895\u00d73=2685'; Replace = '451\u00d78=3608' }\n    @{ Find = '356\u00d75=1780'; Replace = '290\u00d77=2030' }\n    @{ Find = '123\u00d76=738'; Replace = '983\u00d75=4915' }\n    @{ Find = '519\u00d73=1557'; Replace = '530\u00d72=1060' }\n    @{ Find = '768\u00d75=3840'; Replace = '181\u00d79=1629' }\n    @{ Find = '989\u00d74=3956'; Replace = '280\u00d73=840' }\n    @{ Find = '843\u00d77=5901'; Replace = '219\u00d74=876' }\n    @{ Find = '512\u00d77=3584'; Replace = '415\u00d75=2075' }\n    @{ Find = '595\u00d73=1785'; Replace = '553\u00d74=2212' }\n    @{ Find = '571\u00d78=4568'; Replace = '307\u00d75=1535' }\n    @{ Find = '797\u00d77=5579'; Replace = '471\u00d72=942' }\n    @{ Find = '321\u00d79=2889'; Replace = '779\u00d78=6232' }\n    @{ Find = '402\u00d76=2412'; Replace = '271\u00d75=1355' }\n    @{ Find = '222\u00d76=1332'; Replace = '578\u00d75=2890' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $($pair.Find)\"\n    }\n}\n"}
